$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Correct the previously-provisional value for 01-07-2021 (D67: 1.4 -> 1.5)
$ws.Range("D67").Value = 1.5

# Append new monthly row for 01-08-2021
# Force text entry for the date-like label (A68) so it isn't auto-converted
# to a date serial number, then strip the temporary Text format so the
# cell keeps the workbook's default (unstyled) look, matching the rest
# of column A.
$ws.Range("A68").NumberFormat = "@"
$ws.Range("A68").Value = "01-08-2021"
$ws.Range("A68").ClearFormats()

$ws.Range("B68").Value = 0.1
$ws.Range("C68").Value = -0.3
$ws.Range("D68").Value = 0.2
